# mypc_dummy.xlsx - "not perfect mypc finish"
#
# Fills in the "reason" column (D) for the rows whose score (C) is being
# dropped from 100 to 90 - i.e. these students lost 10 points and we now
# record why (PC/password/MAC related excuses instead of the placeholder
# "empty" string).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ten new reason strings, in the exact order they need to be minted
# into the shared-string table (this matches how the workbook author
# originally typed them in, independent of row order). We seed them once
# via a scratch row so the shared-string table gets them in this order,
# then point every target cell at the matching text - Excel's shared
# string dedup takes care of reusing the same entry.
$reasons = @(
  "한컴 업데이트 안됨",
  "패스워드 까먹음",
  "password 까먹음",
  "MAC 노트북 사용 중",
  "맥 노트북 때문에",
  "패스워드 기억안남",
  "한컴 NEO임",
  "비밀번호 까먹음",
  "비밀번호 기억 안남",
  "비번 까먹음"
)

$scratchRow = 1000
for ($i = 0; $i -lt $reasons.Length; $i++) {
  $ws.Cells.Item($scratchRow, $i + 1).Value = $reasons[$i]
}

# row number -> reason text for every row whose score/reason changed
$updates = @{
  9   = "한컴 업데이트 안됨"
  26  = "패스워드 까먹음"
  40  = "한컴 업데이트 안됨"
  67  = "한컴 업데이트 안됨"
  80  = "비밀번호 까먹음"
  89  = "한컴 NEO임"
  95  = "한컴 업데이트 안됨"
  98  = "한컴 업데이트 안됨"
  111 = "한컴 업데이트 안됨"
  112 = "한컴 NEO임"
  119 = "한컴 업데이트 안됨"
  135 = "비밀번호 기억 안남"
  154 = "한컴 NEO임"
  155 = "비번 까먹음"
  177 = "한컴 업데이트 안됨"
  185 = "한컴 업데이트 안됨"
  191 = "한컴 업데이트 안됨"
  192 = "한컴 업데이트 안됨"
  197 = "MAC 노트북 사용 중"
  207 = "한컴 업데이트 안됨"
  208 = "한컴 업데이트 안됨"
  216 = "한컴 업데이트 안됨"
  226 = "한컴 업데이트 안됨"
  232 = "패스워드 기억안남"
  253 = "한컴 업데이트 안됨"
  268 = "한컴 업데이트 안됨"
  277 = "한컴 업데이트 안됨"
  286 = "한컴 업데이트 안됨"
  298 = "한컴 업데이트 안됨"
  311 = "패스워드 까먹음"
  324 = "한컴 업데이트 안됨"
  325 = "맥 노트북 때문에"
  337 = "password 까먹음"
  349 = "한컴 업데이트 안됨"
  358 = "MAC 노트북 사용 중"
}

foreach ($row in $updates.Keys) {
  $ws.Cells.Item($row, 3).Value = 90
  $ws.Cells.Item($row, 4).Value = $updates[$row]
}

# Drop the scratch helper cells again - only used to control shared-string
# creation order, shouldn't remain part of the sheet's used range.
for ($i = 0; $i -lt $reasons.Length; $i++) {
  $ws.Cells.Item($scratchRow, $i + 1).ClearContents()
}

# Leave the UI scrolled/selected where the author ended up.
$ws.Range("H153").Select()
$excel.ActiveWindow.ScrollRow = 139
$excel.ActiveWindow.ScrollColumn = 1
